# Adds 4 new data rows (14-17) to the NIFTY_Options_Analysis sheet.
#
# Row 13 (the last existing data row) is format-copied onto each new row
# first, so every new cell keeps the exact same cellXfs style index as the
# rest of the table (header row untouched). Text columns are briefly marked
# as Text (NumberFormat "@") before the values are typed in, which is what a
# user does via Format Cells > Text to stop Excel silently turning strings
# like "2026-02-20", "10:00:09" or "100%" into a date/time/percent number --
# then row 13 formatting is re-applied on top so the final style matches the
# rest of the sheet precisely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 14 ----
$ws.Range("A13:AE13").Copy()
$ws.Range("A14:AE14").PasteSpecial(-4122)
$ws.Range("A14").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("M14").NumberFormat = "@"
$ws.Range("O14").NumberFormat = "@"
$ws.Range("T14").NumberFormat = "@"
$ws.Range("U14").NumberFormat = "@"
$ws.Range("AC14").NumberFormat = "@"
$ws.Range("AD14").NumberFormat = "@"
$ws.Range("AE14").NumberFormat = "@"

$ws.Range("A14").Value = "2026-02-20"
$ws.Range("B14").Value = "10:00:09"
$ws.Range("C14").Value = "AVOID"
$ws.Range("D14").Value = "AVOID"
$ws.Range("E14").Value = "100%"
$ws.Range("F14").Value = "TRADEABLE"
$ws.Range("G14").Value2 = 0
$ws.Range("H14").Value2 = 25522.2
$ws.Range("I14").Value2 = 14.18
$ws.Range("J14").Value2 = 1.51
$ws.Range("K14").Value2 = 0
$ws.Range("L14").Value2 = 73.09999999999999
$ws.Range("M14").Value = "UNKNOWN"
$ws.Range("N14").Value2 = 0
$ws.Range("O14").Value = "UNKNOWN"
$ws.Range("P14").Value2 = 0
$ws.Range("Q14").Value2 = 0
$ws.Range("R14").Value2 = 0
$ws.Range("S14").Value2 = 0
$ws.Range("T14").Value = "NONE"
$ws.Range("U14").Value = ""
$ws.Range("V14").Value2 = 0
$ws.Range("W14").Value2 = 0
$ws.Range("X14").Value2 = 0
$ws.Range("Y14").Value2 = 0
$ws.Range("Z14").Value2 = 0
$ws.Range("AA14").Value2 = 0
$ws.Range("AB14").Value2 = 0
$ws.Range("AC14").Value = "HARD VETO: CPR TRENDING DAY: Price 25522.20 above TC 25515.24 - BULLISH TRENDING DAY likely"
$ws.Range("AD14").Value = "CPR TRENDING DAY: Price 25522.20 above TC 25515.24 - BULLISH TRENDING DAY likely"
$ws.Range("AE14").Value = "Yes"

$ws.Range("A13:AE13").Copy()
$ws.Range("A14:AE14").PasteSpecial(-4122)

# ---- Row 15 ----
$ws.Range("A13:AE13").Copy()
$ws.Range("A15:AE15").PasteSpecial(-4122)
$ws.Range("A15").NumberFormat = "@"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("M15").NumberFormat = "@"
$ws.Range("O15").NumberFormat = "@"
$ws.Range("T15").NumberFormat = "@"
$ws.Range("U15").NumberFormat = "@"
$ws.Range("AC15").NumberFormat = "@"
$ws.Range("AD15").NumberFormat = "@"
$ws.Range("AE15").NumberFormat = "@"

$ws.Range("A15").Value = "2026-02-23"
$ws.Range("B15").Value = "10:00:12"
$ws.Range("C15").Value = "AVOID"
$ws.Range("D15").Value = "AVOID"
$ws.Range("E15").Value = "100%"
$ws.Range("F15").Value = "TRADEABLE"
$ws.Range("G15").Value2 = 0
$ws.Range("H15").Value2 = 25728.75
$ws.Range("I15").Value2 = 14.23
$ws.Range("J15").Value2 = 2.01
$ws.Range("K15").Value2 = 0
$ws.Range("L15").Value2 = 74.3
$ws.Range("M15").Value = "UNKNOWN"
$ws.Range("N15").Value2 = 0
$ws.Range("O15").Value = "UNKNOWN"
$ws.Range("P15").Value2 = 0
$ws.Range("Q15").Value2 = 0
$ws.Range("R15").Value2 = 0
$ws.Range("S15").Value2 = 0
$ws.Range("T15").Value = "NONE"
$ws.Range("U15").Value = ""
$ws.Range("V15").Value2 = 0
$ws.Range("W15").Value2 = 0
$ws.Range("X15").Value2 = 0
$ws.Range("Y15").Value2 = 0
$ws.Range("Z15").Value2 = 0
$ws.Range("AA15").Value2 = 0
$ws.Range("AB15").Value2 = 0
$ws.Range("AC15").Value = "HARD VETO: CPR TRENDING DAY: Price 25728.75 above TC 25554.72 - BULLISH TRENDING DAY likely"
$ws.Range("AD15").Value = "CPR TRENDING DAY: Price 25728.75 above TC 25554.72 - BULLISH TRENDING DAY likely"
$ws.Range("AE15").Value = "Yes"

$ws.Range("A13:AE13").Copy()
$ws.Range("A15:AE15").PasteSpecial(-4122)

# ---- Row 16 ----
$ws.Range("A13:AE13").Copy()
$ws.Range("A16:AE16").PasteSpecial(-4122)
$ws.Range("A16").NumberFormat = "@"
$ws.Range("B16").NumberFormat = "@"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("M16").NumberFormat = "@"
$ws.Range("O16").NumberFormat = "@"
$ws.Range("T16").NumberFormat = "@"
$ws.Range("U16").NumberFormat = "@"
$ws.Range("AC16").NumberFormat = "@"
$ws.Range("AD16").NumberFormat = "@"
$ws.Range("AE16").NumberFormat = "@"

$ws.Range("A16").Value = "2026-02-24"
$ws.Range("B16").Value = "10:00:12"
$ws.Range("C16").Value = "AVOID"
$ws.Range("D16").Value = "AVOID"
$ws.Range("E16").Value = "100%"
$ws.Range("F16").Value = "TRADEABLE"
$ws.Range("G16").Value2 = 0
$ws.Range("H16").Value2 = 25516.75
$ws.Range("I16").Value2 = 14
$ws.Range("J16").Value2 = 0.54
$ws.Range("K16").Value2 = 0
$ws.Range("L16").Value2 = 71
$ws.Range("M16").Value = "UNKNOWN"
$ws.Range("N16").Value2 = 0
$ws.Range("O16").Value = "UNKNOWN"
$ws.Range("P16").Value2 = 0
$ws.Range("Q16").Value2 = 0
$ws.Range("R16").Value2 = 0
$ws.Range("S16").Value2 = 0
$ws.Range("T16").Value = "NONE"
$ws.Range("U16").Value = ""
$ws.Range("V16").Value2 = 0
$ws.Range("W16").Value2 = 0
$ws.Range("X16").Value2 = 0
$ws.Range("Y16").Value2 = 0
$ws.Range("Z16").Value2 = 0
$ws.Range("AA16").Value2 = 0
$ws.Range("AB16").Value2 = 0
$ws.Range("AC16").Value = "HARD VETO: CPR TRENDING DAY: Price 25516.75 below BC 25690.40 - BEARISH TRENDING DAY likely"
$ws.Range("AD16").Value = "CPR TRENDING DAY: Price 25516.75 below BC 25690.40 - BEARISH TRENDING DAY likely"
$ws.Range("AE16").Value = "Yes"

$ws.Range("A13:AE13").Copy()
$ws.Range("A16:AE16").PasteSpecial(-4122)

# ---- Row 17 ----
$ws.Range("A13:AE13").Copy()
$ws.Range("A17:AE17").PasteSpecial(-4122)
$ws.Range("A17").NumberFormat = "@"
$ws.Range("B17").NumberFormat = "@"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("M17").NumberFormat = "@"
$ws.Range("O17").NumberFormat = "@"
$ws.Range("T17").NumberFormat = "@"
$ws.Range("U17").NumberFormat = "@"
$ws.Range("AC17").NumberFormat = "@"
$ws.Range("AD17").NumberFormat = "@"
$ws.Range("AE17").NumberFormat = "@"

$ws.Range("A17").Value = "2026-02-25"
$ws.Range("B17").Value = "10:00:12"
$ws.Range("C17").Value = "AVOID"
$ws.Range("D17").Value = "AVOID"
$ws.Range("E17").Value = "100%"
$ws.Range("F17").Value = "TRADEABLE"
$ws.Range("G17").Value2 = 0
$ws.Range("H17").Value2 = 25600.15
$ws.Range("I17").Value2 = 13.07
$ws.Range("J17").Value2 = -1.29
$ws.Range("K17").Value2 = 0
$ws.Range("L17").Value2 = 57
$ws.Range("M17").Value = "UNKNOWN"
$ws.Range("N17").Value2 = 0
$ws.Range("O17").Value = "UNKNOWN"
$ws.Range("P17").Value2 = 0
$ws.Range("Q17").Value2 = 0
$ws.Range("R17").Value2 = 0
$ws.Range("S17").Value2 = 0
$ws.Range("T17").Value = "NONE"
$ws.Range("U17").Value = ""
$ws.Range("V17").Value2 = 0
$ws.Range("W17").Value2 = 0
$ws.Range("X17").Value2 = 0
$ws.Range("Y17").Value2 = 0
$ws.Range("Z17").Value2 = 0
$ws.Range("AA17").Value2 = 0
$ws.Range("AB17").Value2 = 0
$ws.Range("AC17").Value = "HARD VETO: CPR TRENDING DAY: Price 25600.15 above TC 25444.67 - BULLISH TRENDING DAY likely"
$ws.Range("AD17").Value = "CPR TRENDING DAY: Price 25600.15 above TC 25444.67 - BULLISH TRENDING DAY likely"
$ws.Range("AE17").Value = "Yes"

$ws.Range("A13:AE13").Copy()
$ws.Range("A17:AE17").PasteSpecial(-4122)

$excel.CutCopyMode = 0
$ws.Range("A1").Select()
